# Highlight the "Watch the Course Overview Video (3min)" bullet (Day 1 -
# Getting Started, Task 1) in yellow, as described in the commit
# "new files added for task 1".
#
# wdYellow = 7 (WdColorIndex)

$d = $word.ActiveDocument

# Locate the paragraph by its visible text so the script is resilient to
# minor paragraph re-numbering elsewhere in the document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Watch the*Course Overview Video*(3min)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Highlight the whole paragraph, including the paragraph mark, so the
    # run(s) AND the pilcrow's rPr (stored on <w:pPr><w:rPr>) get the
    # <w:highlight w:val="yellow"/> entry.
    $target.Range.Font.HighlightColorIndex = 7

    # The hyperlinked run ("Course Overview Video ") lives inside a
    # <w:hyperlink> field; highlighting it via the plain paragraph Range
    # does not reach the field's inner run, so do it explicitly through
    # the Hyperlinks collection.
    $paraStart = $target.Range.Start
    $paraEnd = $target.Range.End
    for ($j = 1; $j -le $d.Hyperlinks.Count; $j++) {
        $hl = $d.Hyperlinks.Item($j)
        if ($hl.Range.Start -ge $paraStart -and $hl.Range.End -le $paraEnd) {
            $hl.Range.Font.HighlightColorIndex = 7
        }
    }
}

Write-Host "Highlighted paragraph: " $target.Range.Text
